$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

# Delete the row containing student Id 8 (Mahmud Eliyev), which is row 6.
# This shifts the remaining row (Id 9, Musa Poladli) up to become the new row 6.
$ws.Rows.Item(6).Delete()
